$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.149.91"
$ws.Range("E2").Value = "  +5.42%  "

# Row 3
$ws.Range("D3").Value = "1.882.12"
$ws.Range("E3").Value = "  +3.98%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "281.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.15%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.01%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5297"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +4.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3538"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.61%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07049"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.97%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.41"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8233"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07820"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.19%  "

# Row 13
$ws.Range("D13").Value = "1.908.03"
$ws.Range("E13").Value = "  +5.46%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.213"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.86"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.72%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9994"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.01%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.62"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.35%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008188"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.30%  "

# Row 19
$ws.Range("E19").Value = "  -0.01%  "

# Row 20
$ws.Range("D20").Value = "27.195.69"
$ws.Range("E20").Value = "  +5.30%  "

# Row 21
$ws.Range("D21").Value = "2.120.68"
$ws.Range("E21").Value = "  +4.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.780"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.16"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.251"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.36%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.411"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +9.83%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.16"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.83%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.65"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.85%  "

# Row 28
$ws.Range("E28").Value = "  +1.31%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.41"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +5.30%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.440"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.51%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.398"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.47%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08966"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.04%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04950"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.96%  "

# Row 34
$ws.Range("E34").Value = "  +4.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7494"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.39%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.905"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.16%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.310"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +9.16%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.427"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.61%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5336"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.68%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01887"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.84%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9736"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.38%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "117.24"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.25%  "

# Row 43
$ws.Range("E43").Value = "  +2.54%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.248"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.16%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4626"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.59%  "

# Row 46
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9987"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1371"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.41%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.459"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.75%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.80"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.79%  "

# Row 50
$ws.Range("E50").Value = "  +2.51%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05958"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.84%  "
